$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update numeric values for rows 2-6 (2014/12 - 2018/12 periods)
$ws.Range("D2").Value = 11250
$ws.Range("E2").Value = 603
$ws.Range("F2").Value = 603
$ws.Range("G2").Value = 647
$ws.Range("H2").Value = 535
$ws.Range("I2").Value = 518
$ws.Range("J2").Value = 17
$ws.Range("K2").Value = 14654
$ws.Range("L2").Value = 9427
$ws.Range("M2").Value = 5227
$ws.Range("N2").Value = 5094
$ws.Range("O2").Value = 132
$ws.Range("P2").Value = 135
$ws.Range("Q2").Value = 330
$ws.Range("R2").Value = -471
$ws.Range("S2").Value = 392
$ws.Range("T2").Value = 464
$ws.Range("U2").Value = -134
$ws.Range("V2").Value = 3470
$ws.Range("W2").Value = 5.36
$ws.Range("X2").Value = 4.75
$ws.Range("AA2").Value = 180.37
$ws.Range("AB2").Value = 3451.01
$ws.Range("AC2").Value = 1926
$ws.Range("AD2").Value = 7.14
$ws.Range("AE2").Value = 18851
$ws.Range("AF2").Value = 0.73
$ws.Range("AG2").Value = 150
$ws.Range("AH2").Value = 1.09
$ws.Range("AI2").Value = 7.83
$ws.Range("AJ2").Value = 27028437
$ws.Range("D3").Value = 22894
$ws.Range("E3").Value = 1185
$ws.Range("F3").Value = 1185
$ws.Range("G3").Value = 1090
$ws.Range("H3").Value = 813
$ws.Range("I3").Value = 796
$ws.Range("J3").Value = 17
$ws.Range("K3").Value = 17040
$ws.Range("L3").Value = 10924
$ws.Range("M3").Value = 6115
$ws.Range("N3").Value = 6010
$ws.Range("O3").Value = 105
$ws.Range("P3").Value = 135
$ws.Range("Q3").Value = 1360
$ws.Range("R3").Value = -2880
$ws.Range("S3").Value = 1272
$ws.Range("T3").Value = 1986
$ws.Range("U3").Value = -625
$ws.Range("V3").Value = 4956
$ws.Range("W3").Value = 5.17
$ws.Range("X3").Value = 3.55
$ws.Range("Y3").Value = 14.34
$ws.Range("Z3").Value = 5.13
$ws.Range("AA3").Value = 178.63
$ws.Range("AB3").Value = 3997.61
$ws.Range("AC3").Value = 2945
$ws.Range("AD3").Value = 5.14
$ws.Range("AE3").Value = 22241
$ws.Range("AF3").Value = 0.68
$ws.Range("AG3").Value = 150
$ws.Range("AH3").Value = 0.99
$ws.Range("AI3").Value = 5.09
$ws.Range("AJ3").Value = 27028437
$ws.Range("D4").Value = 24029
$ws.Range("E4").Value = 1078
$ws.Range("F4").Value = 1078
$ws.Range("G4").Value = 1379
$ws.Range("H4").Value = 939
$ws.Range("I4").Value = 914
$ws.Range("J4").Value = 25
$ws.Range("K4").Value = 18018
$ws.Range("L4").Value = 11222
$ws.Range("M4").Value = 6796
$ws.Range("N4").Value = 6684
$ws.Range("O4").Value = 112
$ws.Range("P4").Value = 135
$ws.Range("Q4").Value = 1165
$ws.Range("R4").Value = -1459
$ws.Range("S4").Value = 261
$ws.Range("T4").Value = 1560
$ws.Range("U4").Value = -395
$ws.Range("V4").Value = 5259
$ws.Range("W4").Value = 4.49
$ws.Range("X4").Value = 3.91
$ws.Range("Y4").Value = 14.4
$ws.Range("Z4").Value = 5.36
$ws.Range("AA4").Value = 165.12
$ws.Range("AB4").Value = 4643.75
$ws.Range("AC4").Value = 3382
$ws.Range("AD4").Value = 4.14
$ws.Range("AE4").Value = 24734
$ws.Range("AF4").Value = 0.57
$ws.Range("AG4").Value = 150
$ws.Range("AH4").Value = 1.07
$ws.Range("AI4").Value = 4.43
$ws.Range("AJ4").Value = 27028437
$ws.Range("D5").Value = 20075
$ws.Range("E5").Value = -122
$ws.Range("F5").Value = -122
$ws.Range("G5").Value = -270
$ws.Range("H5").Value = 304
$ws.Range("I5").Value = 341
$ws.Range("J5").Value = -37
$ws.Range("K5").Value = 17648
$ws.Range("L5").Value = 10922
$ws.Range("M5").Value = 6727
$ws.Range("N5").Value = 6669
$ws.Range("O5").Value = 58
$ws.Range("P5").Value = 135
$ws.Range("Q5").Value = 1409
$ws.Range("R5").Value = -914
$ws.Range("S5").Value = 482
$ws.Range("T5").Value = 1385
$ws.Range("U5").Value = 24
$ws.Range("V5").Value = 5100
$ws.Range("W5").Value = -0.61
$ws.Range("X5").Value = 1.51
$ws.Range("Y5").Value = 5.12
$ws.Range("Z5").Value = 1.7
$ws.Range("AA5").Value = 162.36
$ws.Range("AB5").Value = 4837.96
$ws.Range("AC5").Value = 1263
$ws.Range("AD5").Value = 7.95
$ws.Range("AE5").Value = 24678
$ws.Range("AF5").Value = 0.41
$ws.Range("AG5").Value = 150
$ws.Range("AH5").Value = 1.49
$ws.Range("AI5").Value = 11.87
$ws.Range("AJ5").Value = 27028437
$ws.Range("D6").Value = 19069
$ws.Range("E6").Value = 111
$ws.Range("F6").Value = 111
$ws.Range("G6").Value = 169
$ws.Range("H6").Value = -81
$ws.Range("I6").Value = -66
$ws.Range("K6").Value = 16962
$ws.Range("L6").Value = 10296
$ws.Range("M6").Value = 6666
$ws.Range("N6").Value = 6545
$ws.Range("P6").Value = 135
$ws.Range("Q6").Value = -302
$ws.Range("R6").Value = -2007
$ws.Range("S6").Value = 576
$ws.Range("T6").Value = 1801
$ws.Range("U6").Value = -2103
$ws.Range("V6").Value = 5187
$ws.Range("W6").Value = 0.58
$ws.Range("X6").Value = -0.42
$ws.Range("Y6").Value = -1
$ws.Range("Z6").Value = -0.47
$ws.Range("AA6").Value = 154.45
$ws.Range("AB6").Value = 4727.21
$ws.Range("AC6").Value = -245
$ws.Range("AD6").Value = -18.2
$ws.Range("AE6").Value = 24219
$ws.Range("AF6").Value = 0.18
$ws.Range("AG6").Value = 100
$ws.Range("AH6").Value = 2.25
$ws.Range("AI6").Value = -40.89
$ws.Range("AJ6").Value = 27028437

# Clear cells that should be removed entirely (row 2 extra cols, and rows 7-9 data)
$ws.Range("Y2").ClearContents()
$ws.Range("Z2").ClearContents()
$ws.Range("D7").ClearContents()
$ws.Range("E7").ClearContents()
$ws.Range("G7").ClearContents()
$ws.Range("H7").ClearContents()
$ws.Range("I7").ClearContents()
$ws.Range("K7").ClearContents()
$ws.Range("L7").ClearContents()
$ws.Range("M7").ClearContents()
$ws.Range("N7").ClearContents()
$ws.Range("P7").ClearContents()
$ws.Range("Q7").ClearContents()
$ws.Range("R7").ClearContents()
$ws.Range("S7").ClearContents()
$ws.Range("T7").ClearContents()
$ws.Range("U7").ClearContents()
$ws.Range("W7").ClearContents()
$ws.Range("X7").ClearContents()
$ws.Range("Y7").ClearContents()
$ws.Range("Z7").ClearContents()
$ws.Range("AA7").ClearContents()
$ws.Range("AC7").ClearContents()
$ws.Range("AD7").ClearContents()
$ws.Range("AE7").ClearContents()
$ws.Range("AF7").ClearContents()
$ws.Range("AG7").ClearContents()
$ws.Range("AH7").ClearContents()
$ws.Range("AI7").ClearContents()
$ws.Range("D8").ClearContents()
$ws.Range("E8").ClearContents()
$ws.Range("G8").ClearContents()
$ws.Range("H8").ClearContents()
$ws.Range("I8").ClearContents()
$ws.Range("K8").ClearContents()
$ws.Range("L8").ClearContents()
$ws.Range("M8").ClearContents()
$ws.Range("N8").ClearContents()
$ws.Range("P8").ClearContents()
$ws.Range("Q8").ClearContents()
$ws.Range("R8").ClearContents()
$ws.Range("S8").ClearContents()
$ws.Range("T8").ClearContents()
$ws.Range("U8").ClearContents()
$ws.Range("W8").ClearContents()
$ws.Range("X8").ClearContents()
$ws.Range("Y8").ClearContents()
$ws.Range("Z8").ClearContents()
$ws.Range("AA8").ClearContents()
$ws.Range("AC8").ClearContents()
$ws.Range("AD8").ClearContents()
$ws.Range("AE8").ClearContents()
$ws.Range("AF8").ClearContents()
$ws.Range("AG8").ClearContents()
$ws.Range("AH8").ClearContents()
$ws.Range("AI8").ClearContents()
$ws.Range("D9").ClearContents()
$ws.Range("E9").ClearContents()
$ws.Range("G9").ClearContents()
$ws.Range("H9").ClearContents()
$ws.Range("I9").ClearContents()
$ws.Range("K9").ClearContents()
$ws.Range("L9").ClearContents()
$ws.Range("M9").ClearContents()
$ws.Range("N9").ClearContents()
$ws.Range("P9").ClearContents()
$ws.Range("Q9").ClearContents()
$ws.Range("R9").ClearContents()
$ws.Range("S9").ClearContents()
$ws.Range("T9").ClearContents()
$ws.Range("U9").ClearContents()
$ws.Range("W9").ClearContents()
$ws.Range("X9").ClearContents()
$ws.Range("Y9").ClearContents()
$ws.Range("Z9").ClearContents()
$ws.Range("AA9").ClearContents()
$ws.Range("AC9").ClearContents()
$ws.Range("AD9").ClearContents()
$ws.Range("AE9").ClearContents()
$ws.Range("AF9").ClearContents()
$ws.Range("AG9").ClearContents()
$ws.Range("AH9").ClearContents()
$ws.Range("AI9").ClearContents()

Write-Host "Edit applied"